# Add QMC link to Paul Harris
# Insert a new row in the "Attributes" sheet right after the existing
# j_thomas / app / QMC row (row 56), pushing everything below down by one,
# and populate it with the new p_harris / app / QMC record.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attributes")

# Insert a new row at position 57 (shifts existing rows 57.. down to 58..)
$ws.Rows.Item(57).Insert()

$ws.Cells.Item(57, 1).Value = "p_harris"
$ws.Cells.Item(57, 2).Value = "app"
$ws.Cells.Item(57, 3).Value = "QMC"

# Keep selection / view consistent with the edited state
$ws.Range("A46").Select()
$ws.Range("C55").Select()
